$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 2 (shifting the existing data rows down by 3)
$ws.Rows.Item(2).Resize(3).Insert()

# The insert copies formatting/cells from neighboring rows into columns A, D, E, F
# and applies the header-like style to B:C as well; strip that back down so the
# new rows look like ordinary (unstyled) data rows, matching the rest of the sheet.
$ws.Range("A2:A4").Clear()
$ws.Range("D2:F4").Clear()
$ws.Range("B2:C4").ClearFormats()

# Fill in the new rows with the Loving Caliber songs
$ws.Range("B2").Value = "Loving Caliber - You Set My World On Fire"
$ws.Range("C2").Value = "https://www.youtube.com/watch?v=nQ7SQVXkWr8"

$ws.Range("B3").Value = "Loving Caliber - We Were Dancing In The Dark"
$ws.Range("C3").Value = "https://www.youtube.com/watch?v=P-QYRUPDAQ8"

$ws.Range("B4").Value = "Loving Caliber - I Wish You Were Mine"
$ws.Range("C4").Value = "https://www.youtube.com/watch?v=5j9FKszXLag"
